$wb = $excel.ActiveWorkbook

# --- Sheet "cadastroSucesso" ---
$ws1 = $wb.Worksheets.Item("cadastroSucesso")
$ws1.Range("A2").Value = "BRUN262"

# --- Sheet "CadastroFalha" ---
$ws2 = $wb.Worksheets.Item("CadastroFalha")
$ws2.Range("A2").Value = "BRUN260"
$ws2.Range("C2").Value = "Felipe3"

# Remove the hyperlink attached to B3 while keeping the one on B2.
# The runtime's Hyperlinks.Delete() removes every hyperlink on the sheet,
# so delete them all and re-create only the B2 mail-to link, then restore
# its original (centered) alignment so it keeps using the same cell style.
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:felipe.almeidaa14@gmail.com")
$ws2.Range("B2").HorizontalAlignment = -4108

# Clear row 3 (data + formatting) except column B (keep its style, no value) and column M (stays blank).
$ws2.Range("A3").Clear()
$ws2.Range("C3:L3").Clear()
$ws2.Range("B3").ClearContents()

# Make CadastroFalha the active (selected) sheet/tab with C2 selected.
$ws2.Activate()
$ws2.Range("C2").Select()
